$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.071.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.849.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.50%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.015'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.78%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.014'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '309.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4761'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3685'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.83%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07244'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9334'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.54%  '

# Row 11
$ws.Range("E11").Value = '  +2.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07797'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.62%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.839.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.396'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.54%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.485'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.15%  '

# Row 17
$ws.Range("E17").Value = '  +0.65%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008673'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.11%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.014'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.66%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.121.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.11%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.056'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.08%  '

# Row 23
$ws.Range("E23").Value = '  +0.25%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.940'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.42%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.989'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.79%  '

# Row 28
$ws.Range("E28").Value = '  +0.69%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.926'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08869'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.322'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.75%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.181'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.10%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.521'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.55%  '

# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7379'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.00%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.681'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.64%  '

# Row 36
$ws.Range("E36").Value = '  +3.26%  '

# Row 37
$ws.Range("E37").Value = '  +2.18%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05259'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.04%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.967'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5294'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.53%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.035'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.15%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1526'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.290'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.38%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.57%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4744'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '

# Row 46
$ws.Range("E46").Value = '  +0.65%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.48%  '

# Row 48
$ws.Range("E48").Value = '  +0.79%  '

# Row 49
$ws.Range("E49").Value = '  +2.73%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06061'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.44%  '

# Row 51
$ws.Range("E51").Value = '  +4.02%  '
